# Move the "CofC for sterile" / "CofC for Non-Sterile" rows (A14:B14 and
# A15:B15) further down the sheet to A19:B19 / A20:B20, leaving the C
# column cells (C14/C15) where they are. Rows 16-18 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move A14:B14 -> A19:B19 (values + formatting travel together).
$ws.Range("A14:B14").Cut($ws.Range("A19:B19")) | Out-Null
# Fully clear what's left behind so the emptied cells disappear (no stray
# formatted-but-empty cells), matching the target layout.
$ws.Range("A14:B14").Clear() | Out-Null

# Move A15:B15 -> A20:B20
$ws.Range("A15:B15").Cut($ws.Range("A20:B20")) | Out-Null
$ws.Range("A15:B15").Clear() | Out-Null

# The two source rows no longer carry the tall wrapped-text content, so their
# row height reverts back to the sheet default height.
$ws.Rows(14).AutoFit() | Out-Null
$ws.Rows(15).AutoFit() | Out-Null

# Row 19 now holds the long wrapped text that needs two lines, matching the
# height the row used to have back when it lived at row 14.
$ws.Rows(19).RowHeight = 28.8

# Update the active selection to reflect where editing continued.
$ws.Range("B10").Select() | Out-Null
